# RRHH Liquidacion de sueldos
# Cuando se guarda una liquidacion se hace focus nuevamente en el boton de buscar empleado
#
# This adds a new backlog task row and marks a few existing tasks as
# "terminado" (finished), which causes them to be hidden by the existing
# AutoFilter (which only shows rows with estado = "no comenzado").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark some existing tasks as finished ("terminado"). Row 124 is updated
# after the AutoFilter is re-applied further below so that it is not
# hidden by the filter recalculation (matches the target state, where the
# row keeps showing despite its new status).
$ws.Cells.Item(126, 2).Value = "terminado"
$ws.Cells.Item(127, 2).Value = "terminado"

# Add the new backlog task at the end of the table.
$ws.Cells.Item(129, 1).Value = "rr liquidacion de sueldos, cuando guardo una liquidacion posicionarse en la lupa de buscar empleado"
$ws.Cells.Item(129, 2).Value = "terminado"

# Re-apply the AutoFilter so it covers the extended data range (the filter
# used to stop at row 123, now it must reach row 128) and so hidden rows
# get recomputed for the rows whose status just changed.
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:C128").AutoFilter(2, @("no comenzado"), 7)

# Row 124 switches to "terminado" as well, but stays visible (set after the
# filter so it isn't re-hidden by the recalculation above).
$ws.Cells.Item(124, 2).Value = "terminado"

# Keep the workbook's hidden _FilterDatabase defined name in sync with the
# resized AutoFilter range.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "=Hoja1!`$A`$1:`$C`$128"

# Reflect the cursor position left after the edit.
$null = $ws.Activate()
$null = $ws.Range("B131").Select()
